$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 12.48389999999999

$ws.Range("A4").Value = -21.29150000000002
$ws.Range("C4").Value = -11.33359999999999
$ws.Range("E4").Value = 12.30670000000001

$ws.Range("C5").Value = -14.43640000000001

$ws.Range("A6").Value = -20.17369999999999
$ws.Range("C6").Value = -11.47549999999999

$ws.Range("A7").Value = -21.18360000000002

$ws.Range("A8").Value = -20.55
$ws.Range("C8").Value = -12.2425

$ws.Range("E9").Value = 12.53470000000001

$ws.Range("E11").Value = 13.08219999999999

$ws.Range("E14").Value = 13.77640000000001

$ws.Range("A16").Value = -20.319
$ws.Range("C16").Value = -12.2239

$ws.Range("E18").Value = 13.10709999999999

$ws.Range("A20").Value = -22.83410000000002

$ws.Range("A21").Value = -20.43829999999999

$ws.Range("C22").Value = -11.02379999999999

$ws.Range("E25").Value = 12.97369999999999
